$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.009.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.560.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.07"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.782.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.560.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.021.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0704"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0473"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("E31").Value = "  +3.47%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.421.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("E36").Value = "  +10.17%  "
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0166"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.532"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.697.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("E49").Value = "  +3.35%  "
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0959"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.57%  "
